$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.9642857142857143
$ws.Range("C2").Value = 0.8571428571428571

$ws.Range("B3").Value = 0.9675324675324676
$ws.Range("C3").Value = 0.8733766233766234

$ws.Range("B4").Value = 0.9707792207792207
$ws.Range("C4").Value = 0.8603896103896104

$ws.Range("B5").Value = 0.974025974025974
$ws.Range("C5").Value = 0.8766233766233766

$ws.Range("B6").Value = 0.9675324675324676
$ws.Range("C6").Value = 0.8766233766233766
